# Rename the sheets:
#   "Composite_NL4DV-LLM" -> "NL4DV-LLM Evaluation"
#   "Sheet1"              -> "NL4DV Evaluation"
$wb = $excel.ActiveWorkbook

$wsComposite = $wb.Worksheets.Item("Composite_NL4DV-LLM")
$wsComposite.Name = "NL4DV-LLM Evaluation"

$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Name = "NL4DV Evaluation"

# Make the (renamed) first sheet the active/selected tab, matching the
# tabSelected flag moving from the second sheet to the first.
$wsComposite.Activate()
